$d = $word.ActiveDocument

$pairs = @(
  @("2025-08-05 Tuesday", "2025-08-06 Wednesday"),
  @("347÷4=86, 3", "334÷8=41, 6"),
  @("552÷7=78, 6", "809÷2=404, 1"),
  @("457÷3=152, 1", "663÷2=331, 1"),
  @("115÷4=28, 3", "595÷5=119, 0"),
  @("918÷6=153, 0", "368÷2=184, 0"),
  @("854÷3=284, 2", "823÷3=274, 1"),
  @("336÷5=67, 1", "295÷8=36, 7"),
  @("166÷6=27, 4", "152÷9=16, 8"),
  @("374÷9=41, 5", "108÷7=15, 3"),
  @("430÷7=61, 3", "434÷8=54, 2"),
  @("321÷4=80, 1", "973÷9=108, 1"),
  @("844÷8=105, 4", "702÷8=87, 6"),
  @("424÷7=60, 4", "912÷6=152, 0"),
  @("224÷3=74, 2", "597÷5=119, 2"),
  @("376÷2=188, 0", "646÷6=107, 4"),
  @("176÷5=35, 1", "297÷9=33, 0"),
  @("962÷7=137, 3", "999÷7=142, 5"),
  @("151÷4=37, 3", "847÷3=282, 1"),
  @("379÷3=126, 1", "267÷3=89, 0"),
  @("660÷2=330, 0", "201÷5=40, 1"),
  @("183÷5=36, 3", "254÷5=50, 4"),
  @("985÷9=109, 4", "264÷3=88, 0"),
  @("712÷3=237, 1", "632÷2=316, 0"),
  @("474÷3=158, 0", "763÷8=95, 3"),
  @("748÷4=187, 0", "969÷9=107, 6")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
